$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.000001099251006220214
$ws.Range("C2").Value = 0.3375848360084654
$ws.Range("D2").Value = 2938.103010863317
$ws.Range("E2").Value = 198602002.3250627
$ws.Range("G2").Value = 198604940.7656595
